# Apply keeper selections ("x" marks) to the "Selected keepers" sheet,
# which ripple through the COUNTIFS/SUMIFS formulas on the Summary sheet,
# and restore the Summary tab as the active tab (instead of "Team salary
# worksheet").

$wb = $excel.ActiveWorkbook
$wsKeepers = $wb.Worksheets.Item("Selected keepers")
$wsSummary = $wb.Worksheets.Item("Summary")

# Rows (in the "Selected keepers" sheet) whose column E gets marked "x" -
# these are the players kept by each team for the 2018-19 season.
$xAddr = "E3,E4,E6,E8,E10,E12,E13,E14,E15,E16,E17,E23,E29,E31,E32,E34,E35,E36,E39,E40,E41,E48,E49,E60,E61,E63,E64,E65,E66,E68,E71,E72,E73,E74,E75,E76,E79,E80,E81,E87,E88,E90,E91,E92,E96,E97,E116,E125,E127,E130,E133,E136,E144,E146,E147,E149,E153,E156,E158,E171,E172,E173,E176,E178,E181,E182,E184,E185,E186,E188,E201,E203,E205,E208,E209,E210,E211,E212,E213,E214,E215,E220,E224,E227,E231,E233,E235,E236,E237,E243,E245,E248,E258,E259,E260,E264,E267,E268,E269,E270,E272,E274"
$wsKeepers.Range($xAddr).Value = "x"

# One additional row (33) picks up the same cell formatting as the others
# but is left without a keeper mark (blank "x" cell).
$wsKeepers.Range("E33").Font.Name = "Arial"
$wsKeepers.Range("E33").Font.Size = 10

# Selection bookkeeping on "Selected keepers": cursor parked at A3 (top of
# the frozen pane body) rather than at the default first cell of the pane.
$wsKeepers.Range("A3").Select()

# The workbook now opens on the Summary tab (instead of "Team salary
# worksheet"), with the cursor sitting wherever it was left (A13 clears to
# the sheet default).
$wsSummary.Activate()
$wsSummary.Range("A1").Select()
